$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old footer row (Excel row 188) to accommodate the one net new data row
$ws.Rows.Item(188).Insert()

$n = 186

# Keep the "Conta" (account number) column as text so leading zeros are preserved
$ws.Range("A2:A" + (1+$n)).NumberFormat = "@"

$data = New-Object 'object[,]' $n,3
$data[0,0] = "005305448"
$data[0,1] = "ALPHASITIO"
$data[0,2] = [double]3965175.38
$data[1,0] = "004352384"
$data[1,1] = "BRASFORT"
$data[1,2] = [double]152091.17
$data[2,0] = "004361159"
$data[2,1] = "HFR"
$data[2,2] = [double]59491.2
$data[3,0] = "004212476"
$data[3,1] = "MARIA"
$data[3,2] = [double]50594.84
$data[4,0] = "004480970"
$data[4,1] = "ALBERTO"
$data[4,2] = [double]35255.29
$data[5,0] = "004363260"
$data[5,1] = "LARISSA"
$data[5,2] = [double]28000
$data[6,0] = "005064129"
$data[6,1] = "THIAGO"
$data[6,2] = [double]26552.57
$data[7,0] = "004321016"
$data[7,1] = "JOAQUIM"
$data[7,2] = [double]25206.34
$data[8,0] = "005366255"
$data[8,1] = "RAPHAELA"
$data[8,2] = [double]22526.47
$data[9,0] = "004224011"
$data[9,1] = "THOMAS"
$data[9,2] = [double]16134.66
$data[10,0] = "004690692"
$data[10,1] = "PHYLIA"
$data[10,2] = [double]13007.32
$data[11,0] = "004364200"
$data[11,1] = "BLOCO"
$data[11,2] = [double]10645.98
$data[12,0] = "005274028"
$data[12,1] = "RAFAEL"
$data[12,2] = [double]5051.93
$data[13,0] = "004971783"
$data[13,1] = "ERIKA"
$data[13,2] = [double]5000
$data[14,0] = "005063749"
$data[14,1] = "NATALIA"
$data[14,2] = [double]2330.96
$data[15,0] = "001761119"
$data[15,1] = "BLUEMETRIX"
$data[15,2] = [double]2016.95
$data[16,0] = "004467884"
$data[16,1] = "ANA"
$data[16,2] = [double]1748.17
$data[17,0] = "004854514"
$data[17,1] = "MARCIA"
$data[17,2] = [double]1200.84
$data[18,0] = "004487140"
$data[18,1] = "VALMIR"
$data[18,2] = [double]1160
$data[19,0] = "004504449"
$data[19,1] = "KELMA"
$data[19,2] = [double]1000
$data[20,0] = "004392159"
$data[20,1] = "RODRIGO"
$data[20,2] = [double]900.21
$data[21,0] = "004369172"
$data[21,1] = "LUIZA"
$data[21,2] = [double]679.85
$data[22,0] = "005171652"
$data[22,1] = "BRUNO"
$data[22,2] = [double]291.82
$data[23,0] = "005135532"
$data[23,1] = "FELIPE"
$data[23,2] = [double]280.13
$data[24,0] = "004498637"
$data[24,1] = "TIAGO"
$data[24,2] = [double]243.4
$data[25,0] = "005591536"
$data[25,1] = "GUSTAVO"
$data[25,2] = [double]129.8
$data[26,0] = "004222784"
$data[26,1] = "RAFAEL"
$data[26,2] = [double]125
$data[27,0] = "004958578"
$data[27,1] = "ASSAKO"
$data[27,2] = [double]110.07
$data[28,0] = "005270025"
$data[28,1] = "DENIZE"
$data[28,2] = [double]100
$data[29,0] = "004342617"
$data[29,1] = "JURACI"
$data[29,2] = [double]99.96
$data[30,0] = "005173958"
$data[30,1] = "VENIA"
$data[30,2] = [double]99.86
$data[31,0] = "004754056"
$data[31,1] = "BRUNO"
$data[31,2] = [double]99.6
$data[32,0] = "004472076"
$data[32,1] = "RUBENS"
$data[32,2] = [double]99.18
$data[33,0] = "004339183"
$data[33,1] = "JALISON"
$data[33,2] = [double]95.69
$data[34,0] = "004517506"
$data[34,1] = "LUIZ"
$data[34,2] = [double]94.49
$data[35,0] = "004643880"
$data[35,1] = "GABRIEL"
$data[35,2] = [double]94.09
$data[36,0] = "004431591"
$data[36,1] = "MARIO"
$data[36,2] = [double]93.87
$data[37,0] = "004027477"
$data[37,1] = "GABRIELA"
$data[37,2] = [double]93.28
$data[38,0] = "005055226"
$data[38,1] = "DANILO"
$data[38,2] = [double]89.16
$data[39,0] = "004278033"
$data[39,1] = "DAISY"
$data[39,2] = [double]86.72
$data[40,0] = "004432579"
$data[40,1] = "ANA"
$data[40,2] = [double]86.5
$data[41,0] = "004212132"
$data[41,1] = "JOAO"
$data[41,2] = [double]86.38
$data[42,0] = "004207278"
$data[42,1] = "CESAR"
$data[42,2] = [double]84.93
$data[43,0] = "004239387"
$data[43,1] = "LUIZ"
$data[43,2] = [double]83.93
$data[44,0] = "003115072"
$data[44,1] = "VICTOR"
$data[44,2] = [double]81.27
$data[45,0] = "005348011"
$data[45,1] = "TATIANA"
$data[45,2] = [double]80.91
$data[46,0] = "005009947"
$data[46,1] = "VERANICE"
$data[46,2] = [double]80.79
$data[47,0] = "004318604"
$data[47,1] = "RENAN"
$data[47,2] = [double]80.51
$data[48,0] = "004994036"
$data[48,1] = "BALTASAR"
$data[48,2] = [double]80.5
$data[49,0] = "000330949"
$data[49,1] = "RENATO"
$data[49,2] = [double]80.08
$data[50,0] = "004809902"
$data[50,1] = "PEDRO"
$data[50,2] = [double]79.88
$data[51,0] = "004267976"
$data[51,1] = "E3"
$data[51,2] = [double]79.84
$data[52,0] = "005032151"
$data[52,1] = "ANA"
$data[52,2] = [double]79.27
$data[53,0] = "004454365"
$data[53,1] = "RAFAEL"
$data[53,2] = [double]79.25
$data[54,0] = "005256849"
$data[54,1] = "SANDRO"
$data[54,2] = [double]77.17
$data[55,0] = "004479734"
$data[55,1] = "RODRIGO"
$data[55,2] = [double]76
$data[56,0] = "004453045"
$data[56,1] = "JULIAN"
$data[56,2] = [double]75.83
$data[57,0] = "004230529"
$data[57,1] = "LAIS"
$data[57,2] = [double]75.09
$data[58,0] = "004340984"
$data[58,1] = "RENATA"
$data[58,2] = [double]73.59
$data[59,0] = "004277637"
$data[59,1] = "LARA"
$data[59,2] = [double]73.51
$data[60,0] = "004207374"
$data[60,1] = "ANGELICA"
$data[60,2] = [double]72.95
$data[61,0] = "004520100"
$data[61,1] = "ALEXANDRE"
$data[61,2] = [double]71.72
$data[62,0] = "004268684"
$data[62,1] = "PATRICIA"
$data[62,2] = [double]68.11
$data[63,0] = "004332207"
$data[63,1] = "IRACY"
$data[63,2] = [double]67.97
$data[64,0] = "004452507"
$data[64,1] = "DANIELA"
$data[64,2] = [double]65.08
$data[65,0] = "004756968"
$data[65,1] = "DANIELY"
$data[65,2] = [double]64.89
$data[66,0] = "004855596"
$data[66,1] = "MARIANA"
$data[66,2] = [double]64.36
$data[67,0] = "005186167"
$data[67,1] = "ANDREA"
$data[67,2] = [double]63.64
$data[68,0] = "005305965"
$data[68,1] = "SIDMAR"
$data[68,2] = [double]62.42
$data[69,0] = "004335251"
$data[69,1] = "EDMUNDO"
$data[69,2] = [double]62.39
$data[70,0] = "004242237"
$data[70,1] = "MARIAH"
$data[70,2] = [double]60.32
$data[71,0] = "005348975"
$data[71,1] = "JULIA"
$data[71,2] = [double]60
$data[72,0] = "004321092"
$data[72,1] = "DANIEL"
$data[72,2] = [double]58.87
$data[73,0] = "004451996"
$data[73,1] = "ADRIANO"
$data[73,2] = [double]58.35
$data[74,0] = "004381194"
$data[74,1] = "ALINNE"
$data[74,2] = [double]58.16
$data[75,0] = "004472760"
$data[75,1] = "SANDRA"
$data[75,2] = [double]57.27
$data[76,0] = "004804036"
$data[76,1] = "LUCIANA"
$data[76,2] = [double]56.61
$data[77,0] = "004693308"
$data[77,1] = "LAURA"
$data[77,2] = [double]56.37
$data[78,0] = "004459461"
$data[78,1] = "INTERLAGOS"
$data[78,2] = [double]56.02
$data[79,0] = "004215217"
$data[79,1] = "CAROLINA"
$data[79,2] = [double]55.66
$data[80,0] = "004575632"
$data[80,1] = "ADELE"
$data[80,2] = [double]54.96
$data[81,0] = "005018038"
$data[81,1] = "ELAINE"
$data[81,2] = [double]54.42
$data[82,0] = "005216881"
$data[82,1] = "RENAN"
$data[82,2] = [double]53.5
$data[83,0] = "003836362"
$data[83,1] = "ISABELLA"
$data[83,2] = [double]51.76
$data[84,0] = "004400640"
$data[84,1] = "FELIPE"
$data[84,2] = [double]51.44
$data[85,0] = "004208447"
$data[85,1] = "LEILA"
$data[85,2] = [double]50
$data[86,0] = "005077648"
$data[86,1] = "DUNAS"
$data[86,2] = [double]49.87
$data[87,0] = "004229526"
$data[87,1] = "EDUARDA"
$data[87,2] = [double]49.12
$data[88,0] = "004207184"
$data[88,1] = "CRISTINA"
$data[88,2] = [double]48.93
$data[89,0] = "004466221"
$data[89,1] = "WALTER"
$data[89,2] = [double]48.66
$data[90,0] = "004517080"
$data[90,1] = "TATIANA"
$data[90,2] = [double]47.35
$data[91,0] = "004208733"
$data[91,1] = "REINALDO"
$data[91,2] = [double]46.33
$data[92,0] = "004115403"
$data[92,1] = "HEBERT"
$data[92,2] = [double]45.24
$data[93,0] = "001731007"
$data[93,1] = "GUILHERME"
$data[93,2] = [double]44.59
$data[94,0] = "004272426"
$data[94,1] = "RODRIGO"
$data[94,2] = [double]44.35
$data[95,0] = "005070742"
$data[95,1] = "JUSCELINO"
$data[95,2] = [double]44.06
$data[96,0] = "005000460"
$data[96,1] = "MARIANA"
$data[96,2] = [double]41.03
$data[97,0] = "004397124"
$data[97,1] = "MURYLO"
$data[97,2] = [double]40.43
$data[98,0] = "004752615"
$data[98,1] = "LUZIMAR"
$data[98,2] = [double]39.1
$data[99,0] = "004398174"
$data[99,1] = "DANIELE"
$data[99,2] = [double]39.08
$data[100,0] = "004238164"
$data[100,1] = "DANIELA"
$data[100,2] = [double]38.3
$data[101,0] = "001294033"
$data[101,1] = "VIVIANE"
$data[101,2] = [double]38.26
$data[102,0] = "002401479"
$data[102,1] = "JULIO"
$data[102,2] = [double]37.84
$data[103,0] = "005055239"
$data[103,1] = "NORMAN"
$data[103,2] = [double]37.66
$data[104,0] = "004259649"
$data[104,1] = "BENTO"
$data[104,2] = [double]37.61
$data[105,0] = "004265173"
$data[105,1] = "JULIA"
$data[105,2] = [double]36.93
$data[106,0] = "004973881"
$data[106,1] = "ISABELLA"
$data[106,2] = [double]36.87
$data[107,0] = "004329229"
$data[107,1] = "GABRIEL"
$data[107,2] = [double]36.25
$data[108,0] = "004431689"
$data[108,1] = "LUCA"
$data[108,2] = [double]35.72
$data[109,0] = "004207641"
$data[109,1] = "MAGALI"
$data[109,2] = [double]35.14
$data[110,0] = "004243043"
$data[110,1] = "SUELI"
$data[110,2] = [double]34.15
$data[111,0] = "005437764"
$data[111,1] = "EVA"
$data[111,2] = [double]33.57
$data[112,0] = "004470679"
$data[112,1] = "RODOLFO"
$data[112,2] = [double]33.54
$data[113,0] = "005035754"
$data[113,1] = "JOSE"
$data[113,2] = [double]33.36
$data[114,0] = "005558076"
$data[114,1] = "ALEXANDRE"
$data[114,2] = [double]28.84
$data[115,0] = "004999434"
$data[115,1] = "EDUARDO"
$data[115,2] = [double]28.06
$data[116,0] = "004377415"
$data[116,1] = "ANGELA"
$data[116,2] = [double]26.37
$data[117,0] = "004940560"
$data[117,1] = "CRISTIANO"
$data[117,2] = [double]25.34
$data[118,0] = "004240292"
$data[118,1] = "MARCO"
$data[118,2] = [double]24.3
$data[119,0] = "004404724"
$data[119,1] = "LEANDRO"
$data[119,2] = [double]24.14
$data[120,0] = "005581299"
$data[120,1] = "ZILDA"
$data[120,2] = [double]23.99
$data[121,0] = "005009922"
$data[121,1] = "ANA"
$data[121,2] = [double]23.62
$data[122,0] = "004259659"
$data[122,1] = "BENTO"
$data[122,2] = [double]22.99
$data[123,0] = "004563237"
$data[123,1] = "FERNANDO"
$data[123,2] = [double]21.74
$data[124,0] = "005219257"
$data[124,1] = "CAROLINE"
$data[124,2] = [double]21.37
$data[125,0] = "004214604"
$data[125,1] = "MARIA"
$data[125,2] = [double]20.75
$data[126,0] = "004204255"
$data[126,1] = "AMADO"
$data[126,2] = [double]18.77
$data[127,0] = "004527606"
$data[127,1] = "MARCIA"
$data[127,2] = [double]18.68
$data[128,0] = "004181486"
$data[128,1] = "ANDREA"
$data[128,2] = [double]18
$data[129,0] = "004479463"
$data[129,1] = "HENRIQUE"
$data[129,2] = [double]17.36
$data[130,0] = "004946997"
$data[130,1] = "EDUARDO"
$data[130,2] = [double]16.98
$data[131,0] = "005143579"
$data[131,1] = "GABRIEL"
$data[131,2] = [double]16.18
$data[132,0] = "005169333"
$data[132,1] = "EDUARDO"
$data[132,2] = [double]16.12
$data[133,0] = "004920447"
$data[133,1] = "MARILIA"
$data[133,2] = [double]16.02
$data[134,0] = "004995535"
$data[134,1] = "ASIEL"
$data[134,2] = [double]15.76
$data[135,0] = "005133039"
$data[135,1] = "PAULO"
$data[135,2] = [double]15.48
$data[136,0] = "004848279"
$data[136,1] = "LEONARDO"
$data[136,2] = [double]15.25
$data[137,0] = "005135281"
$data[137,1] = "RAFAEL"
$data[137,2] = [double]15.2
$data[138,0] = "004581652"
$data[138,1] = "CINCO"
$data[138,2] = [double]14.94
$data[139,0] = "004422594"
$data[139,1] = "WANDIR"
$data[139,2] = [double]14.67
$data[140,0] = "004472538"
$data[140,1] = "RODOLFO"
$data[140,2] = [double]14.25
$data[141,0] = "002064834"
$data[141,1] = "RAFAELA"
$data[141,2] = [double]13.24
$data[142,0] = "004207955"
$data[142,1] = "SILVANIA"
$data[142,2] = [double]12.48
$data[143,0] = "005366671"
$data[143,1] = "TATIANA"
$data[143,2] = [double]11.45
$data[144,0] = "004752461"
$data[144,1] = "SERGIO"
$data[144,2] = [double]10.77
$data[145,0] = "004216298"
$data[145,1] = "FLORDELIZ"
$data[145,2] = [double]9.76
$data[146,0] = "002894447"
$data[146,1] = "JOAO"
$data[146,2] = [double]9.17
$data[147,0] = "004264780"
$data[147,1] = "MARCELO"
$data[147,2] = [double]8.99
$data[148,0] = "004479965"
$data[148,1] = "DIEGO"
$data[148,2] = [double]8.78
$data[149,0] = "004693631"
$data[149,1] = "NELY"
$data[149,2] = [double]7.36
$data[150,0] = "004767746"
$data[150,1] = "ISABELE"
$data[150,2] = [double]7.33
$data[151,0] = "004419765"
$data[151,1] = "WALTER"
$data[151,2] = [double]7.21
$data[152,0] = "004508159"
$data[152,1] = "FELIPE"
$data[152,2] = [double]7.11
$data[153,0] = "004530494"
$data[153,1] = "ROSANGELA"
$data[153,2] = [double]6.94
$data[154,0] = "004854496"
$data[154,1] = "JOSE"
$data[154,2] = [double]6.64
$data[155,0] = "004289402"
$data[155,1] = "LARISSA"
$data[155,2] = [double]6.47
$data[156,0] = "004805133"
$data[156,1] = "PATRICIA"
$data[156,2] = [double]6.22
$data[157,0] = "004448501"
$data[157,1] = "JOAO"
$data[157,2] = [double]5.55
$data[158,0] = "005142624"
$data[158,1] = "RODRIGO"
$data[158,2] = [double]4.75
$data[159,0] = "004462543"
$data[159,1] = "RODOLFO"
$data[159,2] = [double]4.3
$data[160,0] = "005242683"
$data[160,1] = "LUCAS"
$data[160,2] = [double]4.26
$data[161,0] = "004848927"
$data[161,1] = "ULDARICO"
$data[161,2] = [double]3.62
$data[162,0] = "004936634"
$data[162,1] = "LEONARDO"
$data[162,2] = [double]3.08
$data[163,0] = "004308815"
$data[163,1] = "ZELI"
$data[163,2] = [double]1.25
$data[164,0] = "004431546"
$data[164,1] = "GABRIELA"
$data[164,2] = [double]0.97
$data[165,0] = "004223502"
$data[165,1] = "BRUNA"
$data[165,2] = [double]0.78
$data[166,0] = "005624274"
$data[166,1] = "CLAYTON"
$data[166,2] = [double]0.68
$data[167,0] = "002694089"
$data[167,1] = "VITOR"
$data[167,2] = [double]0.65
$data[168,0] = "004453302"
$data[168,1] = "ISABELLA"
$data[168,2] = [double]0.39
$data[169,0] = "004806286"
$data[169,1] = "VERA"
$data[169,2] = [double]0.19
$data[170,0] = "004371857"
$data[170,1] = "NAZARETH"
$data[170,2] = [double]0.18
$data[171,0] = "004357159"
$data[171,1] = "JOAO"
$data[171,2] = [double]0.15
$data[172,0] = "004320840"
$data[172,1] = "NATALIA"
$data[172,2] = [double]0.14
$data[173,0] = "004335031"
$data[173,1] = "EDMUNDO"
$data[173,2] = [double]0.11
$data[174,0] = "004466350"
$data[174,1] = "RAQUEL"
$data[174,2] = [double]0.11
$data[175,0] = "005047946"
$data[175,1] = "GABRIEL"
$data[175,2] = [double]0.09
$data[176,0] = "004589311"
$data[176,1] = "CLARICE"
$data[176,2] = [double]0.06
$data[177,0] = "001090818"
$data[177,1] = "MARIA"
$data[177,2] = [double]0.01
$data[178,0] = "002878817"
$data[178,1] = "GUILHERME"
$data[178,2] = [double]0.01
$data[179,0] = "004384258"
$data[179,1] = "PAULA"
$data[179,2] = [double]0.01
$data[180,0] = "004400000"
$data[180,1] = "VILMA"
$data[180,2] = [double]0.01
$data[181,0] = "004474935"
$data[181,1] = "MELISSA"
$data[181,2] = [double]0.01
$data[182,0] = "004475136"
$data[182,1] = "MATHEO"
$data[182,2] = [double]0.01
$data[183,0] = "004612043"
$data[183,1] = "YURI"
$data[183,2] = [double]0.01
$data[184,0] = "004974089"
$data[184,1] = "CELIA"
$data[184,2] = [double]0.01
$data[185,0] = "005314853"
$data[185,1] = "GLEUBER"
$data[185,2] = [double]0.01

$rng = $ws.Range("A2:C" + (1+$n))
$rng.Value = $data

Write-Host ("Data written: " + $n + " rows; total used rows now: " + $ws.UsedRange.Rows.Count)
